$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 43.666668
$ws.Range("I12").Value = 45.5
$ws.Range("K12").Value = 45.5
$ws.Range("M12").Value = 124.5

$ws.Range("H15").Value = 1170.6666
$ws.Range("I15").Value = 1170.6666
$ws.Range("K15").Value = 3511.9998
$ws.Range("M15").Value = -3342.9998

$ws.Range("H28").Value = 611
$ws.Range("I28").Value = 514.6667
$ws.Range("K28").Value = 514.6667
$ws.Range("M28").Value = -29.66669999999999

$ws.Range("H43").Value = 6800
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 6800
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 6800
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -6938

$ws.Range("H92").Value = 1049.5
$ws.Range("I92").Value = 1049.5
$ws.Range("K92").Value = 1049.5
$ws.Range("M92").Value = 198.5

$ws.Range("H121").Value = 1589.75
$ws.Range("J121").Value = 1589.75
$ws.Range("L121").Value = 4769.25
$ws.Range("N121").Value = -8263.25

$ws.Range("H137").Value = 4970.1177
$ws.Range("I137").Value = 2497.25
$ws.Range("J137").Value = 5731
$ws.Range("K137").Value = 7491.75
$ws.Range("L137").Value = 17193
$ws.Range("M137").Value = -4941.75
$ws.Range("N137").Value = -22293

$ws.Range("H138").Value = 5217.3335
$ws.Range("I138").Value = 2000
$ws.Range("J138").Value = 5378.2
$ws.Range("K138").Value = 6000
$ws.Range("L138").Value = 16134.6
$ws.Range("M138").Value = -860
$ws.Range("N138").Value = -26414.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 850
$ws.Range("J4").Value = 925
$ws.Range("L4").Value = 925
$ws.Range("N4").Value = -1157

$ws.Range("H32").Value = 2638.6287
$ws.Range("I32").Value = 1995.5758
$ws.Range("J32").Value = 13249
$ws.Range("K32").Value = 1995.5758
$ws.Range("L32").Value = 13249
$ws.Range("M32").Value = -1708.5758
$ws.Range("N32").Value = -13823

$ws.Range("H45").Value = 1067.75
$ws.Range("I45").Value = 985.36365
$ws.Range("K45").Value = 985.36365
$ws.Range("M45").Value = -608.36365

$ws.Range("H61").Value = 3090.25
$ws.Range("I61").Value = 3098.7273
$ws.Range("J61").Value = 2997
$ws.Range("K61").Value = 3098.7273
$ws.Range("L61").Value = 2997
$ws.Range("M61").Value = -2886.7273
$ws.Range("N61").Value = -3421

$ws.Range("H74").Value = 1000
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -126
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1000
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -632
$ws.Range("N77").ClearContents()

$ws.Range("H122").Value = 3994.4
$ws.Range("I122").Value = 4093.9092
$ws.Range("J122").Value = 3720.75
$ws.Range("K122").Value = 12281.7276
$ws.Range("L122").Value = 11162.25
$ws.Range("M122").Value = -9831.7276
$ws.Range("N122").Value = -16062.25

$ws.Range("H132").Value = 4441.3335
$ws.Range("I132").Value = 2994.8
$ws.Range("K132").Value = 8984.400000000001
$ws.Range("M132").Value = -6454.400000000001

$ws.Range("H136").Value = 3090.25
$ws.Range("I136").Value = 3098.7273
$ws.Range("J136").Value = 2997
$ws.Range("K136").Value = 9296.1819
$ws.Range("L136").Value = 8991
$ws.Range("M136").Value = -6746.1819
$ws.Range("N136").Value = -14091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 18001.75
$ws.Range("I20").Value = 17666.166
$ws.Range("K20").Value = 17666.166
$ws.Range("M20").Value = -17419.166

$ws.Range("H86").Value = 539.3333
$ws.Range("I86").Value = 447.2
$ws.Range("K86").Value = 447.2
$ws.Range("M86").Value = 675.8

$ws.Range("H89").Value = 539.3333
$ws.Range("I89").Value = 447.2
$ws.Range("K89").Value = 2236
$ws.Range("M89").Value = 3380

$ws.Range("H107").Value = 3751.8333
$ws.Range("I107").Value = 3751.8333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3751.8333
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1831.8333
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 2224.7334
$ws.Range("I134").Value = 2224.7334
$ws.Range("K134").Value = 6674.2002
$ws.Range("M134").Value = -4139.2002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 423.5
$ws.Range("J7").Value = 999
$ws.Range("L7").Value = 999
$ws.Range("N7").Value = -1225

$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H31").Value = 10201.833
$ws.Range("I31").Value = 3304.75
$ws.Range("J31").Value = 23996
$ws.Range("K31").Value = 3304.75
$ws.Range("L31").Value = 23996
$ws.Range("M31").Value = -3009.75
$ws.Range("N31").Value = -24586

$ws.Range("H34").Value = 10201.833
$ws.Range("I34").Value = 3304.75
$ws.Range("J34").Value = 23996
$ws.Range("K34").Value = 3304.75
$ws.Range("L34").Value = 23996
$ws.Range("M34").Value = -3102.75
$ws.Range("N34").Value = -24400

$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H132").Value = 5165.3335
$ws.Range("I132").Value = 3499.5
$ws.Range("J132").Value = 5998.25
$ws.Range("K132").Value = 10498.5
$ws.Range("L132").Value = 17994.75
$ws.Range("M132").Value = -7968.5
$ws.Range("N132").Value = -23054.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 44.25
$ws.Range("I2").Value = 49
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 294
$ws.Range("L2").Value = 180
$ws.Range("M2").Value = -181
$ws.Range("N2").Value = -406

$ws.Range("H34").Value = 2477.6
$ws.Range("J34").Value = 2972
$ws.Range("L34").Value = 8916
$ws.Range("N34").Value = -9084

$ws.Range("H38").Value = 101
$ws.Range("J38").Value = 157
$ws.Range("L38").Value = 471
$ws.Range("N38").Value = -1165

$ws.Range("H55").Value = 2499.5
$ws.Range("J55").Value = 2499.5
$ws.Range("L55").Value = 7498.5
$ws.Range("N55").Value = -7852.5

$ws.Range("H80").Value = 9002
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 9002
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2157.111
$ws.Range("I102").Value = 2130.7144
$ws.Range("K102").Value = 2130.7144
$ws.Range("M102").Value = -508.7143999999998

$ws.Range("H136").Value = 133333.33
$ws.Range("J136").Value = 133333.33
$ws.Range("L136").Value = 399999.99
$ws.Range("N136").Value = -405099.99

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 943.25
$ws.Range("J16").Value = 494.66666
$ws.Range("L16").Value = 494.66666
$ws.Range("N16").Value = -834.66666

$ws.Range("H40").Value = 4500
$ws.Range("I40").Value = 4500
$ws.Range("K40").Value = 4500
$ws.Range("M40").Value = -4364

$ws.Range("H61").Value = 1416.5
$ws.Range("I61").Value = 500
$ws.Range("K61").Value = 500
$ws.Range("M61").Value = -298

$ws.Range("H113").Value = 1416.5
$ws.Range("I113").Value = 500
$ws.Range("K113").Value = 500
$ws.Range("M113").Value = 1670

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3795.3333
$ws.Range("I81").Value = 1769.5
$ws.Range("K81").Value = 3539
$ws.Range("M81").Value = -2478

$ws.Range("H84").Value = 3795.3333
$ws.Range("I84").Value = 1769.5
$ws.Range("K84").Value = 17695
$ws.Range("M84").Value = -12391

$ws.Range("H107").Value = 356.8
$ws.Range("I107").Value = 341.1111
$ws.Range("K107").Value = 1023.3333
$ws.Range("M107").Value = 896.6667

$ws.Range("H113").Value = 440.1
$ws.Range("I113").Value = 466.57144
$ws.Range("J113").Value = 378.33334
$ws.Range("K113").Value = 1399.71432
$ws.Range("L113").Value = 1135.00002
$ws.Range("M113").Value = 770.28568
$ws.Range("N113").Value = -5475.000019999999

$ws.Range("H122").Value = 1566
$ws.Range("I122").Value = 1566
$ws.Range("K122").Value = 4698
$ws.Range("M122").Value = -2248

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 8826.223
$ws.Range("I136").Value = 8463.058999999999
$ws.Range("K136").Value = 25389.177
$ws.Range("M136").Value = -22839.177

$ws.Range("H137").Value = 44996.332
$ws.Range("J137").Value = 44996.332
$ws.Range("L137").Value = 44996.332
$ws.Range("N137").Value = -55196.332
